# Apply "Add data for 2022-09-05" update to carjacking-by-month-yoy-latest.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab and the workbook's sheet entry (date moves from 08-27 to 08-28)
$ws.Name = "Through 2022-08-28"

# Update the "August (through 08-27)" label cell to "August (through 08-28)"
$ws.Range("A9").Value = "August (through 08-28)"

# Update August row (row 9) values for columns C..I
$ws.Cells.Item(9, 3).Value = 68
$ws.Cells.Item(9, 4).Value = 82
$ws.Cells.Item(9, 5).Value = 59
$ws.Cells.Item(9, 6).Value = 42
$ws.Cells.Item(9, 7).Value = 154
$ws.Cells.Item(9, 8).Value = 149
$ws.Cells.Item(9, 9).Value = 149

# Update Total row (row 10) values for columns C..I
$ws.Cells.Item(10, 3).Value = 370
$ws.Cells.Item(10, 4).Value = 547
$ws.Cells.Item(10, 5).Value = 484
$ws.Cells.Item(10, 6).Value = 346
$ws.Cells.Item(10, 7).Value = 775
$ws.Cells.Item(10, 8).Value = 1059
$ws.Cells.Item(10, 9).Value = 1120
